$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.395158389027614
$ws.Range("C2").Value = 0.05243438464364658
$ws.Range("D2").Value = 0.0784572293712813
$ws.Range("E2").Value = 0.4117279779247696
$ws.Range("G2").Value = 0.002448044911749249
$ws.Range("K2").Value = 0.3709183826733238
$ws.Range("O2").Value = 3.496825343056088

$ws.Range("B3").Value = 0.3548139464674875
$ws.Range("C3").Value = 0.0466696856260711
$ws.Range("D3").Value = 0.07117407217988614
$ws.Range("E3").Value = 0.3592284641779031
$ws.Range("G3").Value = 0.002451319617835405
$ws.Range("K3").Value = 0.32747339459209
$ws.Range("O3").Value = 3.480233362765802

$ws.Range("B4").Value = 0.3301482408787138
$ws.Range("C4").Value = 0.04310844417855719
$ws.Range("D4").Value = 0.06673819535259895
$ws.Range("E4").Value = 0.3270895470033679
$ws.Range("G4").Value = 0.002453434886117125
$ws.Range("K4").Value = 0.3008400178826207
$ws.Range("O4").Value = 3.472161266396057

$ws.Range("B5").Value = 0.3201235615563576
$ws.Range("C5").Value = 0.04165175845837155
$ws.Range("D5").Value = 0.06493955397233719
$ws.Range("E5").Value = 0.314014947791776
$ws.Range("G5").Value = 0.002454323260730585
$ws.Range("K5").Value = 0.2899974001738315
$ws.Range("O5").Value = 3.469402149964623

$ws.Range("B6").Value = 0.3184605943464476
$ws.Range("C6").Value = 0.0414095478307388
$ws.Range("D6").Value = 0.06464143421638369
$ws.Range("E6").Value = 0.3118452112761361
$ws.Range("G6").Value = 0.002454472370672374
$ws.Range("K6").Value = 0.288197642099135
$ws.Range("O6").Value = 3.468975984146027

$ws.Range("B7").Value = 0.3300129358414665
$ws.Range("C7").Value = 0.04308882086723997
$ws.Range("D7").Value = 0.0667139017998295
$ws.Range("E7").Value = 0.3269131306605715
$ws.Range("G7").Value = 0.002453446760035132
$ws.Range("K7").Value = 0.3006937470435105
$ws.Range("O7").Value = 3.472121910875615

$ws.Range("B8").Value = 0.3812257146403795
$ws.Range("C8").Value = 0.05045121916515427
$ws.Range("D8").Value = 0.07593849563612309
$ws.Range("E8").Value = 0.3936050789105394
$ws.Range("G8").Value = 0.002449152376195946
$ws.Range("K8").Value = 0.3559298932890158
$ws.Range("O8").Value = 3.490664377782196

$ws.Range("B9").Value = 0.4824932362587049
$ws.Range("C9").Value = 0.06471729413848948
$ws.Range("D9").Value = 0.09431646179486108
$ws.Range("E9").Value = 0.525239996082405
$ws.Range("G9").Value = 0.002441556965165385
$ws.Range("K9").Value = 0.4645812857680482
$ws.Range("O9").Value = 3.543890712459842

$ws.Range("B10").Value = 0.5574128951548118
$ws.Range("C10").Value = 0.0750958234608845
$ws.Range("D10").Value = 0.1080000219128294
$ws.Range("E10").Value = 0.6226060295658158
$ws.Range("G10").Value = 0.002436474503841463
$ws.Range("K10").Value = 0.5446195630510999
$ws.Range("O10").Value = 3.593396698710563

$ws.Range("B11").Value = 0.5916104020959381
$ws.Range("C11").Value = 0.07979545041146707
$ws.Range("D11").Value = 0.1142656305106584
$ws.Range("E11").Value = 0.6670720617880761
$ws.Range("G11").Value = 0.002434269280767354
$ws.Range("K11").Value = 0.5810796406980785
$ws.Range("O11").Value = 3.618203123497153

$ws.Range("B12").Value = 0.6045768124897961
$ws.Range("C12").Value = 0.08157198613147898
$ws.Range("D12").Value = 0.1166442025419059
$ws.Range("E12").Value = 0.6839374234294411
$ws.Range("G12").Value = 0.002433449489186747
$ws.Range("K12").Value = 0.5948933821596825
$ws.Range("O12").Value = 3.627927367939265

$ws.Range("B13").Value = 0.6017835299568901
$ws.Range("C13").Value = 0.08118951576417999
$ws.Range("D13").Value = 0.1161316704343562
$ws.Range("E13").Value = 0.6803039304099201
$ws.Range("G13").Value = 0.002433625367855815
$ws.Range("K13").Value = 0.5919180318723534
$ws.Range("O13").Value = 3.625818344766458

$ws.Range("B14").Value = 0.5926768254943227
$ws.Range("C14").Value = 0.07994166943264247
$ws.Range("D14").Value = 0.1144611981324459
$ws.Range("E14").Value = 0.6684590314433905
$ws.Range("G14").Value = 0.002434201530328778
$ws.Range("K14").Value = 0.5822159638968856
$ws.Range("O14").Value = 3.618996506618544

$ws.Range("B15").Value = 0.5871008579090642
$ws.Range("C15").Value = 0.07917692233095863
$ws.Range("D15").Value = 0.113438758101978
$ws.Range("E15").Value = 0.6612072717196185
$ws.Range("G15").Value = 0.002434556433831766
$ws.Range("K15").Value = 0.5762740883048707
$ws.Range("O15").Value = 3.614861042319546

$ws.Range("B16").Value = 0.5551803325106732
$ws.Range("C16").Value = 0.07478825654366972
$ws.Range("D16").Value = 0.107591376043672
$ws.Range("E16").Value = 0.6197037411215689
$ws.Range("G16").Value = 0.002436620762689669
$ws.Range("K16").Value = 0.5422378116710149
$ws.Range("O16").Value = 3.591821699115826

$ws.Range("B17").Value = 0.5356277715143847
$ws.Range("C17").Value = 0.07209041321847565
$ws.Range("D17").Value = 0.1040146990343089
$ws.Range("E17").Value = 0.594288623638505
$ws.Range("G17").Value = 0.002437914461591997
$ws.Range("K17").Value = 0.5213704396381047
$ws.Range("O17").Value = 3.578274616441689

$ws.Range("B18").Value = 0.5243925993429741
$ws.Range("C18").Value = 0.07053665035547851
$ws.Range("D18").Value = 0.1019613350803752
$ws.Range("E18").Value = 0.5796867487102588
$ws.Range("G18").Value = 0.00243866862175088
$ws.Range("K18").Value = 0.5093728173144427
$ws.Range("O18").Value = 3.570697739599154

$ws.Range("B19").Value = 0.5205904473655778
$ws.Range("C19").Value = 0.0700102235725808
$ws.Range("D19").Value = 0.1012667606097182
$ws.Range("E19").Value = 0.5747455296733648
$ws.Range("G19").Value = 0.002438925697385486
$ws.Range("K19").Value = 0.505311446307303
$ws.Range("O19").Value = 3.56816921675329

$ws.Range("B20").Value = 0.5377080434783181
$ws.Range("C20").Value = 0.07237781385563835
$ws.Range("D20").Value = 0.1043950440520121
$ws.Range("E20").Value = 0.5969924118926144
$ws.Range("G20").Value = 0.002437775704633511
$ws.Range("K20").Value = 0.5235913175779672
$ws.Range("O20").Value = 3.579694456378547

$ws.Range("B21").Value = 0.5953512371588658
$ws.Range("C21").Value = 0.08030827659240458
$ws.Range("D21").Value = 0.1149516952789327
$ws.Range("E21").Value = 0.6719374150171404
$ws.Range("G21").Value = 0.002434031883439696
$ws.Range("K21").Value = 0.5850655046841098
$ws.Range("O21").Value = 3.620991258690538

$ws.Range("B22").Value = 0.6331208370745003
$ws.Range("C22").Value = 0.08547315894580265
$ws.Range("D22").Value = 0.1218856152535608
$ws.Range("E22").Value = 0.7210768901219922
$ws.Range("G22").Value = 0.002431674098427937
$ws.Range("K22").Value = 0.6252837786418297
$ws.Range("O22").Value = 3.649908917893526

$ws.Range("B23").Value = 0.6129537140849664
$ws.Range("C23").Value = 0.08271822367396453
$ws.Range("D23").Value = 0.1181816773156896
$ws.Range("E23").Value = 0.6948350311105713
$ws.Range("G23").Value = 0.002432924374119306
$ws.Range("K23").Value = 0.603814786150906
$ws.Range("O23").Value = 3.634297999684179

$ws.Range("B24").Value = 0.5367675342756115
$ws.Range("C24").Value = 0.07224788855418751
$ws.Range("D24").Value = 0.10422308099389
$ws.Range("E24").Value = 0.5957699994738732
$ws.Range("G24").Value = 0.002437838404445892
$ws.Range("K24").Value = 0.5225872608125144
$ws.Range("O24").Value = 3.579051887960276

$ws.Range("B25").Value = 0.4550071581536486
$ws.Range("C25").Value = 0.06087612532209619
$ws.Range("D25").Value = 0.08931329060094129
$ws.Range("E25").Value = 0.489524219498648
$ws.Range("G25").Value = 0.002443523889318158
$ws.Range("K25").Value = 0.4351516234246446
$ws.Range("O25").Value = 3.527673798270342

